# Apply update: "Add data for 2022-03-26"
# - Rename the "through" date from 2022-03-17 to 2022-03-18 (sheet tab name + header cell)
# - Update several cell counts across the neighborhood-by-month grid

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet tab and update the column header text for the "through" date
$ws.Name = "Through 2022-03-18"
$ws.Range("B1").Value = "March 2022 (through March 18)"

# Row 3 - Austin
$ws.Range("H3").Value = 3

# Row 4 - North Lawndale
$ws.Range("T4").Value = 3

# Row 5 - Garfield Park
$ws.Range("H5").Value = 4

# Row 9 - Chicago Lawn
$ws.Range("T9").Value = 1

# Row 11 - Englewood
$ws.Range("B11").Value = 7
$ws.Range("H11").Value = 3

# Row 14 - West Town
$ws.Range("B14").Value = 1

# Row 16 - Little Italy, UIC
$ws.Range("W16").Value = 3

# Row 17 - Auburn Gresham
$ws.Range("W17").Value = 1

# Row 26 - Grand Crossing
$ws.Range("T26").Value = 2

# Row 31 - Near South Side
$ws.Range("B31").Value = 2

# Row 32 - New City
$ws.Range("B32").Value = 3
$ws.Range("W32").Value = 1

# Row 34 - River North
$ws.Range("K34").Value = 1

# Row 36 - Roseland
$ws.Range("T36").Value = 2

# Row 41 - Loop
$ws.Range("B41").Value = 3

# Row 44 - Grand Boulevard
$ws.Range("E44").Value = 2

# Row 57 - Douglas
$ws.Range("E57").Value = 1

# Row 67 - Hyde Park
$ws.Range("Q67").Value = 1

# Row 80 - Riverdale
$ws.Range("N80").Value = 3
